$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text blocks reused further down (copied verbatim from the workbook so the
# shared-string content is byte-identical to the source data).
$semestral        = 'Semestral'
$docente          = '5840494 - Maria Eleonora Andrade de Carvalho'
$shortSyllabusEn  = 'Introduction to metabolism. Bioenergetics. Biological oxidations. Glycides transport  metabolism,  photosynthesis,  lipid metabolism, metabolism of nitrogenous compounds: amino acids, integration and control of metabolic processes, vital cycles: oxygen, carbon, nitrogen and sulfur.'
$longSyllabusEn   = 'Introduction to metabolism. Proteins, polysaccharides, lipids.: catabolic and anabolic pathways. Bioenergetics. The direction of processes: free energy: reaction with the balance constant and with redox potential. Exergonic processes. Role of phosphate: potential of transference of phosphate group. Role of ATP as a free energy currency. Biological oxidations. Electron and proton transporters coenzymes: nucleotides, flavin nucleotides, coenzyme Q. Pyridine and flavin nucleotides- dehydrogenase dependent. Oxidases.Structure of mitochondrial membrane. Respiratory chain: function. Oxidative phosphorylation. Transport. Composition of biological membranes: chemical constitution, characterization, barrier, permeability. Carriers and canal ionophores. Transportation: mediated and non-mediated. Glycides  metabolism. Anaerobic and aerobic degradation of glycides: glycolysis  localization of enzymes, reactions, NaDH.H+ production, the first substrt level phosphorylation, energetic balance; Citric acid cycle  localization of operating enzymes, reactions, production of reduced coenzymes, energetic balance. the pentose phosphate pathway: localization of the enzymes, reactions (oxidative and non-oxidative phases), NaDH.H+ production (physiological implication). Fermentation: definition, fermentation and respiration, raw-materials used in the starch and sugar fermentations, metabolic fates of pyruvate: ethanol and lactic metabolism, acetic and citric. Photosynthesis. Chloroplasts structure. The light reactions. Role of the chlorophyll in the photosynthesis. Cyclic and non-cyclic phosphorylation. NADP reduction. Water photolysis. Synthesis of the acceptor of CO2, Ru-1, 5-diP. Calvin cycle. Lipids  metabolism. Beta-oxidation pathways, oxidation of fatty acids with odd-numbered carbon chains, energetic balance of beta-oxidation, alpha-oxidation and w-oxidation.  Glycerol metabolism. Formation of ketone bodies. Biosynthesis of fatty acids. Amino acids  catabolic pathways. Digestion of proteins, aspects of amino acid synthesis and degradation.Transamination, urea cycle. Metabolic integration. Common metabolites to the metabolism of glycides, lipids and amino acids. Vital cycles: oxygen, carbon, nitrogen and sulfur. Carbon and oxygen cycles. Nitrogen cycle: biological fixation, nitrification, use of nitrate, incorporation of ammonia in organic compounds. Sulfur cycle: sulfate assimilation.'
$dataAtivacao     = '01/01/2018'
$avaliacaoMetodo  = 'A avaliação será feita por meio de provas escritas.'
$notaFinal        = 'A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2'
$recuperacao      = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$requisitoFraco   = 'LOT2007 -  Bioquímica I  (Requisito fraco)
'

# Row 10 ("Objetivos:") now shows the responsible-professor line instead of the
# long Portuguese objectives paragraph.
$ws.Range("B10").Value = $docente
$ws.Range("C10").Value = $docente

# Row 13 ("Programa resumido:") drops its long paragraph in favour of "Semestral"
# and gains the label in column A (it previously had none).
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = $semestral
$ws.Range("C13").Value = $semestral

# Row 14 becomes "Short syllabus:" with the English short-syllabus paragraph.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = $shortSyllabusEn
$ws.Range("C14").Value = $shortSyllabusEn

# Row 15 becomes "Programa:" but its data column now just repeats the activation date.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = $dataAtivacao
$ws.Range("C15").Value = $dataAtivacao

# Row 16 becomes "Syllabus:" with the long English syllabus paragraph.
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = $longSyllabusEn
$ws.Range("C16").Value = $longSyllabusEn

# Row 17 becomes "Avaliação:" and no longer carries any B/C content.
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# Row 18 becomes "Método:" and now shows the responsible-professor line.
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = $docente
$ws.Range("C18").Value = $docente

# Row 19 becomes "Critério:" with the evaluation-method sentence.
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = $avaliacaoMetodo
$ws.Range("C19").Value = $avaliacaoMetodo

# Row 20 becomes "Norma de recuperação:" with the final-grade formula.
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = $notaFinal
$ws.Range("C20").Value = $notaFinal

# Row 21 becomes "Bibliografia:" with the recovery-exam paragraph.
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao

# Row 22 becomes "Requisitos:" and drops the bibliography paragraph entirely.
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# Row 23 drops its "Requisitos:" label (moved up to row 22) and instead carries the
# weak-prerequisite text that used to live on row 24.
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = $requisitoFraco
$ws.Range("C23").Value = $requisitoFraco

# The sheet now only has 23 rows of data, so delete the old trailing row 24.
$ws.Rows.Item(24).Delete()

# --- Row heights ---
# Row 13 gains the 60pt custom height used by the other single-paragraph rows.
$ws.Rows.Item(13).RowHeight = 60

# Row 15 grows from 60pt to 120pt (it now sits where the long-paragraph rows are).
$ws.Rows.Item(15).RowHeight = 120

# Row 17 no longer needs a custom height (used to be 120pt).
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(17).AutoFit()

# Row 18 gains the 60pt height used by the other single-paragraph rows.
$ws.Rows.Item(18).RowHeight = 60

# Row 21 grows from 60pt to 120pt.
$ws.Rows.Item(21).RowHeight = 120

# Row 23 keeps the 30pt height that used to belong to row 24.
$ws.Rows.Item(23).RowHeight = 30
